# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets.
# 展览: F2 142 -> 144, F3 37 -> 38
# 全部类型: F2 142 -> 144, F3 37 -> 38

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 144
    $ws.Range("F3").Value = 38
}
